# saving back dismantle year
# - bump the maximum investment capacity per year (MW) cap way up
# - clarify its comment text
# - drop the unused scratch columns E:F on the "Coupling Parameters" sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# maximum_investment_capacity_per_year: 15000 -> 1000000
$ws.Range("B24").Value = 1000000

# Update its accompanying note text
$ws.Range("C24").Value = "MW. Planned power plants from the input are also considered. So the maximum should be large"

# Remove the leftover scratch/testing data in columns E:F
$ws.Columns("E:F").Delete() | Out-Null

# Reflect where the user's selection ended up
$ws.Activate() | Out-Null
$ws.Range("C14").Select() | Out-Null
